$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-12 Friday", 2) | Out-Null
$d.Content.Find.Execute("640÷6=106, 4", $true, $false, $false, $false, $false, $true, 1, $false, "373÷5=74, 3", 2) | Out-Null
$d.Content.Find.Execute("300÷2=150, 0", $true, $false, $false, $false, $false, $true, 1, $false, "995÷9=110, 5", 2) | Out-Null
$d.Content.Find.Execute("419÷6=69, 5", $true, $false, $false, $false, $false, $true, 1, $false, "276÷4=69, 0", 2) | Out-Null
$d.Content.Find.Execute("620÷8=77, 4", $true, $false, $false, $false, $false, $true, 1, $false, "176÷3=58, 2", 2) | Out-Null
$d.Content.Find.Execute("496÷6=82, 4", $true, $false, $false, $false, $false, $true, 1, $false, "408÷2=204, 0", 2) | Out-Null
$d.Content.Find.Execute("855÷7=122, 1", $true, $false, $false, $false, $false, $true, 1, $false, "116÷2=58, 0", 2) | Out-Null
$d.Content.Find.Execute("195÷6=32, 3", $true, $false, $false, $false, $false, $true, 1, $false, "869÷7=124, 1", 2) | Out-Null
$d.Content.Find.Execute("961÷3=320, 1", $true, $false, $false, $false, $false, $true, 1, $false, "441÷2=220, 1", 2) | Out-Null
$d.Content.Find.Execute("915÷9=101, 6", $true, $false, $false, $false, $false, $true, 1, $false, "631÷5=126, 1", 2) | Out-Null
$d.Content.Find.Execute("626÷5=125, 1", $true, $false, $false, $false, $false, $true, 1, $false, "619÷4=154, 3", 2) | Out-Null
$d.Content.Find.Execute("472÷9=52, 4", $true, $false, $false, $false, $false, $true, 1, $false, "826÷7=118, 0", 2) | Out-Null
$d.Content.Find.Execute("719÷6=119, 5", $true, $false, $false, $false, $false, $true, 1, $false, "164÷8=20, 4", 2) | Out-Null
$d.Content.Find.Execute("894÷9=99, 3", $true, $false, $false, $false, $false, $true, 1, $false, "453÷3=151, 0", 2) | Out-Null
$d.Content.Find.Execute("938÷3=312, 2", $true, $false, $false, $false, $false, $true, 1, $false, "985÷6=164, 1", 2) | Out-Null
$d.Content.Find.Execute("757÷4=189, 1", $true, $false, $false, $false, $false, $true, 1, $false, "619÷4=154, 3", 2) | Out-Null
$d.Content.Find.Execute("606÷6=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "315÷7=45, 0", 2) | Out-Null
$d.Content.Find.Execute("579÷5=115, 4", $true, $false, $false, $false, $false, $true, 1, $false, "163÷4=40, 3", 2) | Out-Null
$d.Content.Find.Execute("215÷7=30, 5", $true, $false, $false, $false, $false, $true, 1, $false, "768÷8=96, 0", 2) | Out-Null
$d.Content.Find.Execute("129÷3=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "748÷7=106, 6", 2) | Out-Null
$d.Content.Find.Execute("342÷2=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "124÷7=17, 5", 2) | Out-Null
$d.Content.Find.Execute("107÷7=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "146÷8=18, 2", 2) | Out-Null
$d.Content.Find.Execute("465÷6=77, 3", $true, $false, $false, $false, $false, $true, 1, $false, "720÷7=102, 6", 2) | Out-Null
$d.Content.Find.Execute("634÷4=158, 2", $true, $false, $false, $false, $false, $true, 1, $false, "953÷4=238, 1", 2) | Out-Null
$d.Content.Find.Execute("221÷7=31, 4", $true, $false, $false, $false, $false, $true, 1, $false, "666÷9=74, 0", 2) | Out-Null
$d.Content.Find.Execute("555÷4=138, 3", $true, $false, $false, $false, $false, $true, 1, $false, "567÷2=283, 1", 2) | Out-Null
